# "update scripts wuth new tpm" -- NATMI re-ran with new TPM values.
# ECs is no longer a "Sending cluster" (only FAPs/MuSCs send now), so the
# 3 rows where Sending cluster = ECs (old rows 8-10) are removed, and the
# remaining 6 rows (old rows 2-7, FAPs/MuSCs sending) get refreshed
# edge-expression numbers for the Gdf1->Bmpr1a ligand-receptor pair.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused rows 8-10 (ECs as sending cluster is dropped)
$ws.Rows("8:10").Delete()

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Gdf1"
$ws.Range("C2").Value = "Bmpr1a"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.006993666666666666
$ws.Range("H2").Value = 0.020981
$ws.Range("I2").Value = 0.4853567132414176
$ws.Range("J2").Value = 0.4853567132414175
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.443038
$ws.Range("N2").Value = 4.329114
$ws.Range("O2").Value = 0.0289666880885598
$ws.Range("P2").Value = 0.0289666880885598
$ws.Range("Q2").Value = 0.01009212675933333
$ws.Range("R2").Value = 0.090829140834
$ws.Range("S2").Value = 0.01405917652415271
$ws.Range("T2").Value = 0.0140591765241527

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Gdf1"
$ws.Range("C3").Value = "Bmpr1a"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.006993666666666666
$ws.Range("H3").Value = 0.020981
$ws.Range("I3").Value = 0.4853567132414176
$ws.Range("J3").Value = 0.4853567132414175
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 29.20351433333333
$ws.Range("N3").Value = 87.610543
$ws.Range("O3").Value = 0.5862140087672342
$ws.Range("P3").Value = 0.5862140087672342
$ws.Range("Q3").Value = 0.2042396447425555
$ws.Range("R3").Value = 1.838156802683
$ws.Range("S3").Value = 0.2845229045513403
$ws.Range("T3").Value = 0.2845229045513403

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Gdf1"
$ws.Range("C4").Value = "Bmpr1a"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.006993666666666666
$ws.Range("H4").Value = 0.020981
$ws.Range("I4").Value = 0.4853567132414176
$ws.Range("J4").Value = 0.4853567132414175
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 19.170603
$ws.Range("N4").Value = 57.511809
$ws.Range("O4").Value = 0.384819303144206
$ws.Range("P4").Value = 0.384819303144206
$ws.Range("Q4").Value = 0.134072807181
$ws.Range("R4").Value = 1.206655264629
$ws.Range("S4").Value = 0.1867746321659245
$ws.Range("T4").Value = 0.1867746321659245

# Row 5
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Gdf1"
$ws.Range("C5").Value = "Bmpr1a"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.007415666666666667
$ws.Range("H5").Value = 0.022247
$ws.Range("I5").Value = 0.5146432867585824
$ws.Range("J5").Value = 0.5146432867585824
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.443038
$ws.Range("N5").Value = 4.329114
$ws.Range("O5").Value = 0.0289666880885598
$ws.Range("P5").Value = 0.0289666880885598
$ws.Range("Q5").Value = 0.01070108879533333
$ws.Range("R5").Value = 0.09630979915799999
$ws.Range("S5").Value = 0.01490751156440709
$ws.Range("T5").Value = 0.01490751156440709

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Gdf1"
$ws.Range("C6").Value = "Bmpr1a"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.007415666666666667
$ws.Range("H6").Value = 0.022247
$ws.Range("I6").Value = 0.5146432867585824
$ws.Range("J6").Value = 0.5146432867585824
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 29.20351433333333
$ws.Range("N6").Value = 87.610543
$ws.Range("O6").Value = 0.5862140087672342
$ws.Range("P6").Value = 0.5862140087672342
$ws.Range("Q6").Value = 0.2165635277912222
$ws.Range("R6").Value = 1.949071750121
$ws.Range("S6").Value = 0.3016911042158938
$ws.Range("T6").Value = 0.3016911042158938

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Gdf1"
$ws.Range("C7").Value = "Bmpr1a"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.007415666666666667
$ws.Range("H7").Value = 0.022247
$ws.Range("I7").Value = 0.5146432867585824
$ws.Range("J7").Value = 0.5146432867585824
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 19.170603
$ws.Range("N7").Value = 57.511809
$ws.Range("O7").Value = 0.384819303144206
$ws.Range("P7").Value = 0.384819303144206
$ws.Range("Q7").Value = 0.142162801647
$ws.Range("R7").Value = 1.279465214823
$ws.Range("S7").Value = 0.1980446709782815
$ws.Range("T7").Value = 0.1980446709782814
